$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 2.45
$ws.Range("I2").Value = 3.5
$ws.Range("J2").Value = 3.4
$ws.Range("Q2").Value = 2.38
$ws.Range("R2").Value = 1.59
$ws.Range("U2").Value = 5.8
$ws.Range("V2").Value = 1.14
$ws.Range("AD2").Value = 9.5
$ws.Range("AE2").Value = 11
$ws.Range("AF2").Value = 23
$ws.Range("AN2").Value = 15
$ws.Range("AO2").Value = 15
# Row 3
$ws.Range("G3").Value = 2.25
$ws.Range("I3").Value = 3.9
$ws.Range("J3").Value = 3.1
$ws.Range("L3").Value = 4.5
$ws.Range("M3").Value = 1.13
$ws.Range("N3").Value = 6
$ws.Range("AD3").Value = 9
$ws.Range("AF3").Value = 21
$ws.Range("AG3").Value = 23
$ws.Range("AJ3").Value = 5.5
$ws.Range("AM3").Value = 8
# Row 4
$ws.Range("G4").Value = 3.25
$ws.Range("I4").Value = 2.4
$ws.Range("J4").Value = 4.33
$ws.Range("K4").Value = 1.8
$ws.Range("L4").Value = 3.4
$ws.Range("M4").Value = 1.14
$ws.Range("N4").Value = 5.5
$ws.Range("AA4").Value = 2.38
$ws.Range("AB4").Value = 1.53
$ws.Range("AE4").Value = 13
$ws.Range("AF4").Value = 41
$ws.Range("AI4").Value = 5.5
$ws.Range("AK4").Value = 23
$ws.Range("AL4").Value = 101
$ws.Range("AM4").Value = 5.5
$ws.Range("AN4").Value = 9.5
$ws.Range("AP4").Value = 23
# Row 5
$ws.Range("U5").Value = 4.9
$ws.Range("V5").Value = 1.18
# Row 6
$ws.Range("G6").Value = 1.38
$ws.Range("H6").Value = 5
$ws.Range("I6").Value = 7.5
$ws.Range("J6").Value = 1.91
$ws.Range("L6").Value = 8.5
$ws.Range("Q6").Value = 1.49
$ws.Range("R6").Value = 2.65
$ws.Range("S6").Value = 2
$ws.Range("T6").Value = 1.85
$ws.Range("U6").Value = 2.75
$ws.Range("V6").Value = 1.46
$ws.Range("AA6").Value = 2.25
$ws.Range("AB6").Value = 1.57
$ws.Range("AC6").Value = 5.5
$ws.Range("AD6").Value = 5.5
$ws.Range("AF6").Value = 8.5
$ws.Range("AI6").Value = 10
$ws.Range("AL6").Value = 101
$ws.Range("AM6").Value = 15
$ws.Range("AN6").Value = 41
$ws.Range("AO6").Value = 23
$ws.Range("AP6").Value = 101
$ws.Range("AQ6").Value = 67
$ws.Range("AR6").Value = 67
# Row 8
$ws.Range("G8").Value = 2.1
$ws.Range("I8").Value = 3.25
$ws.Range("K8").Value = 2.1
$ws.Range("S8").Value = 1.98
$ws.Range("T8").Value = 1.88
$ws.Range("AC8").Value = 8
# Row 9
$ws.Range("G9").Value = 1.6
$ws.Range("I9").Value = 5.25
$ws.Range("K9").Value = 2.4
$ws.Range("L9").Value = 5
$ws.Range("M9").Value = 1.03
$ws.Range("N9").Value = 15
$ws.Range("O9").Value = 1.18
$ws.Range("P9").Value = 4.5
$ws.Range("S9").Value = 1.6
$ws.Range("T9").Value = 2.3
$ws.Range("U9").Value = 1.98
$ws.Range("V9").Value = 1.83
$ws.Range("W9").Value = 2.5
$ws.Range("X9").Value = 1.5
$ws.Range("Y9").Value = 1.3
$ws.Range("Z9").Value = 3.4
$ws.Range("AA9").Value = 1.67
$ws.Range("AB9").Value = 2.1
$ws.Range("AC9").Value = 9
$ws.Range("AD9").Value = 9
$ws.Range("AF9").Value = 13
$ws.Range("AH9").Value = 21
$ws.Range("AI9").Value = 15
$ws.Range("AS9").Value = 151
# Row 13
$ws.Range("G13").Value = 3.5
$ws.Range("H13").Value = 3.3
$ws.Range("I13").Value = 2.15
$ws.Range("L13").Value = 2.88
$ws.Range("M13").Value = 1.07
$ws.Range("N13").Value = 9
$ws.Range("AC13").Value = 9.5
$ws.Range("AE13").Value = 12
$ws.Range("AN13").Value = 10
$ws.Range("AQ13").Value = 19
# Row 14
$ws.Range("G14").Value = 3.9
$ws.Range("H14").Value = 3.6
$ws.Range("J14").Value = 4.5
$ws.Range("L14").Value = 2.6
$ws.Range("M14").Value = 1.06
$ws.Range("N14").Value = 10
$ws.Range("O14").Value = 1.3
$ws.Range("P14").Value = 3.4
$ws.Range("S14").Value = 2
$ws.Range("T14").Value = 1.85
$ws.Range("Y14").Value = 1.4
$ws.Range("Z14").Value = 2.75
$ws.Range("AA14").Value = 1.8
$ws.Range("AB14").Value = 1.95
$ws.Range("AG14").Value = 34
$ws.Range("AH14").Value = 41
$ws.Range("AI14").Value = 10
$ws.Range("AM14").Value = 7
$ws.Range("AN14").Value = 9
$ws.Range("AS14").Value = 301
# Row 15
$ws.Range("G15").Value = 2.27
$ws.Range("H15").Value = 3.85
$ws.Range("I15").Value = 2.62
$ws.Range("J15").Value = 2.65
$ws.Range("K15").Value = 2.52
$ws.Range("L15").Value = 3
$ws.Range("O15").Value = 1.11
$ws.Range("P15").Value = 5.6
$ws.Range("S15").Value = 1.35
$ws.Range("T15").Value = 2.95
$ws.Range("W15").Value = 1.85
$ws.Range("X15").Value = 1.85
$ws.Range("Y15").Value = 1.21
$ws.Range("Z15").Value = 4
$ws.Range("AA15").Value = 1.33
$ws.Range("AB15").Value = 3.05
$ws.Range("AC15").Value = 17
$ws.Range("AF15").Value = 28
$ws.Range("AG15").Value = 15.5
$ws.Range("AH15").Value = 16
$ws.Range("AI15").Value = 10.25
$ws.Range("AJ15").Value = 9
$ws.Range("AK15").Value = 10.5
$ws.Range("AL15").Value = 25
$ws.Range("AM15").Value = 17
$ws.Range("AN15").Value = 20
$ws.Range("AO15").Value = 10.5
$ws.Range("AP15").Value = 35
$ws.Range("AQ15").Value = 18
$ws.Range("AS15").Value = 110
# Row 16
$ws.Range("H16").Value = 4.05
$ws.Range("K16").Value = 2.35
$ws.Range("M16").Value = 1.03
$ws.Range("N16").Value = 9
$ws.Range("P16").Value = 4.2
$ws.Range("AC16").Value = 17.5
$ws.Range("AI16").Value = 9
$ws.Range("AJ16").Value = 8.25
$ws.Range("AN16").Value = 8.5
$ws.Range("AQ16").Value = 11.5
# Row 17
$ws.Range("H17").Value = 2.8
$ws.Range("I17").Value = 3.25
$ws.Range("L17").Value = 4
$ws.Range("M17").Value = 1.13
$ws.Range("N17").Value = 6
$ws.Range("AA17").Value = 2.1
$ws.Range("AB17").Value = 1.63
$ws.Range("AI17").Value = 6
$ws.Range("AM17").Value = 7.5
$ws.Range("AO17").Value = 13

$wb.Save()